# Aggiunto supporto per tab diverse per apertura e chiusura valvola.
#
# 1) Rinomina il foglio "Misurazioni" in "MisurazioniApertura"
#    (Excel aggiorna automaticamente i riferimenti nei Nomi definiti
#    "Angolo" e "Trimmer" che puntano a quel foglio).
# 2) Aggiunge un nuovo foglio "MisurazioniChiusura" subito dopo
#    "MisurazioniApertura" (cioe' prima di "Grafico").
# 3) Ripristina "DatiGenerali" come foglio attivo/selezionato, cosi'
#    come lo era prima della modifica.

$wb = $excel.ActiveWorkbook

$apertura = $wb.Worksheets.Item("Misurazioni")
$apertura.Name = "MisurazioniApertura"

$grafico = $wb.Worksheets.Item("Grafico")
$chiusura = $wb.Worksheets.Add($grafico)
$chiusura.Name = "MisurazioniChiusura"

$datiGenerali = $wb.Worksheets.Item("DatiGenerali")
$datiGenerali.Activate()
